$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a literal TEXT value (never let Excel auto-convert
# numeric-looking strings like "214.41" into a number), while leaving the
# cell style exactly as it was (no lingering "@" text format).
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = '29.907.21'
$ws.Range("E2").Value = '  +0.49%  '

# Row 3
$ws.Range("D3").Value = '1.632.80'
$ws.Range("E3").Value = '  +1.84%  '

# Row 4
$ws.Range("E4").Value = '  +0.30%  '

# Row 5
Set-TextValue $ws.Range("D5") '214.41'
$ws.Range("E5").Value = '  +0.88%  '

# Row 6
$ws.Range("E6").Value = '  +0.31%  '

# Row 7
$ws.Range("E7").Value = '  +0.28%  '

# Row 8
Set-TextValue $ws.Range("D8") '28.43'
$ws.Range("E8").Value = '  +0.79%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.257'
$ws.Range("E9").Value = '  +1.36%  '

# Row 10
$ws.Range("E10").Value = '  +0.66%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.0911'
$ws.Range("E11").Value = '  +0.25%  '

# Row 12
$ws.Range("D12").Value = '1.866.61'
$ws.Range("E12").Value = '  +1.86%  '

# Row 13
$ws.Range("D13").Value = '1.634.59'
$ws.Range("E13").Value = '  +1.99%  '

# Row 14
$ws.Range("E14").Value = '  +2.63%  '

# Row 15
$ws.Range("E15").Value = '  +17.75%  '

# Row 16
$ws.Range("D16").Value = '29.947.96'
$ws.Range("E16").Value = '  +0.63%  '

# Row 17
$ws.Range("E17").Value = '  +2.18%  '

# Row 18
Set-TextValue $ws.Range("D18") '63.96'
$ws.Range("E18").Value = '  -0.33%  '

# Row 19
Set-TextValue $ws.Range("D19") '242.05'
$ws.Range("E19").Value = '  +0.08%  '

# Row 20
$ws.Range("E20").Value = '  +0.44%  '

# Row 21
$ws.Range("E21").Value = '  +0.21%  '

# Row 22
Set-TextValue $ws.Range("D22") '9.83'
$ws.Range("E22").Value = '  +4.50%  '

# Row 23
$ws.Range("E23").Value = '  +2.41%  '

# Row 24
$ws.Range("E24").Value = '  +1.52%  '

# Row 25
Set-TextValue $ws.Range("D25") '157.75'
$ws.Range("E25").Value = '  +1.64%  '

# Row 26
Set-TextValue $ws.Range("D26") '15.49'
$ws.Range("E26").Value = '  +0.31%  '

# Row 27
$ws.Range("E27").Value = '  +0.19%  '

# Row 29
$ws.Range("E29").Value = '  +0.25%  '

# Row 30
$ws.Range("E30").Value = '  +1.72%  '

# Row 31
$ws.Range("E31").Value = '  +4.18%  '

# Row 32
$ws.Range("E32").Value = '  +3.98%  '

# Row 33
$ws.Range("E33").Value = '  -0.59%  '

# Row 34
$ws.Range("D34").Value = '1.425.45'
$ws.Range("E34").Value = '  +0.29%  '

# Row 35
$ws.Range("E35").Value = '  +4.81%  '

# Row 36
$ws.Range("E36").Value = '  -0.33%  '

# Row 37
$ws.Range("E37").Value = '  -3.80%  '

# Row 38
$ws.Range("E38").Value = '  +0.04%  '

# Row 39
$ws.Range("E39").Value = '  +0.56%  '

# Row 40
Set-TextValue $ws.Range("D40") '75.87'
$ws.Range("E40").Value = '  +12.69%  '

# Row 41
$ws.Range("E41").Value = '  +1.23%  '

# Row 42
$ws.Range("E42").Value = '  +2.34%  '

# Row 43
$ws.Range("E43").Value = '  +1.08%  '

# Row 44
$ws.Range("E44").Value = '  -1.50%  '

# Row 45
$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue $ws.Range("D45") '1.00'
$ws.Range("E45").Value = '  +0.23%  '

# Row 46
$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range("D46") '1.02'
$ws.Range("E46").Value = '  +3.27%  '

# Row 47
$ws.Range("B47").Value = 'BitcoinSV'
$ws.Range("C47").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextValue $ws.Range("D47") '52.94'
$ws.Range("E47").Value = '  -5.98%  '

# Row 48
$ws.Range("D48").Value = '1.774.91'
$ws.Range("E48").Value = '  +2.05%  '

# Row 49
Set-TextValue $ws.Range("D49") '5.35'
$ws.Range("E49").Value = '  -0.58%  '

# Row 50
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Range("D50") '90.10'
$ws.Range("E50").Value = '  +4.10%  '

# Row 51
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₆0112'
$ws.Range("E51").Value = '  +8.40%  '
